$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace header labels "Comp.N" with "PCN" for columns A1:L1
$ws.Range("A1").Value = "PC1"
$ws.Range("B1").Value = "PC2"
$ws.Range("C1").Value = "PC3"
$ws.Range("D1").Value = "PC4"
$ws.Range("E1").Value = "PC5"
$ws.Range("F1").Value = "PC6"
$ws.Range("G1").Value = "PC7"
$ws.Range("H1").Value = "PC8"
$ws.Range("I1").Value = "PC9"
$ws.Range("J1").Value = "PC10"
$ws.Range("K1").Value = "PC11"
$ws.Range("L1").Value = "PC12"
